# Apply the two textual edits described by the diff:
#   1. "pulling at the edges, " + "&" (two separate runs) -> merged into a
#      single run "pulling at the edges, &" (same run formatting, just
#      collapsed into one <w:r>).
#   2. " unmoving, &" -> " motionless, &" (plain word substitution).

$d = $word.ActiveDocument

# --- Edit 1: merge the "pulling at the edges, " / "&" runs -------------
# Re-using the exact same text as the replacement forces Word to rewrite
# the matched span as a single run, collapsing the run boundary that used
# to sit between "edges, " and "&".
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(
    "pulling at the edges, &", $true, $false, $false, $false, $false,
    $true, 1, $false, "pulling at the edges, &", 2)
Write-Host "Merge runs found: $found1"

# --- Edit 2: " unmoving, &" -> " motionless, &" -------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    " unmoving, &", $true, $false, $false, $false, $false,
    $true, 1, $false, " motionless, &", 2)
Write-Host "Word replace found: $found2"
